# New weekly price record for "Terminal La Palmera de La Serena - Arándano (blue)".
# The existing data rows (2021 data, sorted by entry order) are pushed down by one
# row and a new row is inserted at row 8 holding the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 8:15 down to 9:16, opening up a blank row 8.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with this week's record.
$ws.Cells.Item(8, 1).Value2 = 8
$ws.Cells.Item(8, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(8, 3).Value2 = "Coquimbo"
$ws.Cells.Item(8, 4).Value2 = 44494
$ws.Cells.Item(8, 5).Value2 = 4
$ws.Cells.Item(8, 6).Value2 = "Fruta"
$ws.Cells.Item(8, 7).Value2 = 100101
$ws.Cells.Item(8, 8).Value2 = "Berries"
$ws.Cells.Item(8, 9).Value2 = 100101001
$ws.Cells.Item(8, 10).Value2 = "Arándano (blue)"
$ws.Cells.Item(8, 11).Value2 = "Sin especificar"
$ws.Cells.Item(8, 12).Value2 = "Primera"
$ws.Cells.Item(8, 13).Value2 = 200
$ws.Cells.Item(8, 14).Value2 = 11500
$ws.Cells.Item(8, 15).Value2 = 12000
$ws.Cells.Item(8, 16).Value2 = 11750
$ws.Cells.Item(8, 17).Value2 = "`$/bandeja 2 kilos"
$ws.Cells.Item(8, 18).Value2 = "Provincia de Limarí"
$ws.Cells.Item(8, 19).Value2 = 5875
$ws.Cells.Item(8, 20).Value2 = 2
